# "fix big lương phụ phẫu 2 đơn thu nợ 1 tháng. Chưa fix được"
#
# The workbook is a Notion export ("Lũy kế tháng CẦN THƠ"). Several rows
# (pages) share the same Notion last_edited_time string in column D; that
# shared timestamp advanced, and the numeric metrics on the affected page's
# row (row 7) were corrected upward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế tháng CẦN THƠ")

# All cells that carried the old last_edited_time "2024-08-24T20:33:00.000Z"
# now show the new timestamp "2024-08-26T17:26:00.000Z".
$lastEditedCells = @("D4", "D5", "D6", "D7", "D8", "D12", "D13")
foreach ($addr in $lastEditedCells) {
    $ws.Range($addr).Value = "2024-08-26T17:26:00.000Z"
}

# Numeric corrections on row 7 (the page that was fixed).
$ws.Range("T7").Value = 35500000
$ws.Range("W7").Value = 99810000
$ws.Range("AA7").Value = 115040000
$ws.Range("AE7").Value = 214850000
$ws.Range("AH7").Value = 164850000
$ws.Range("AK7").Value = 32
$ws.Range("AN7").Value = 50000000
$ws.Range("AQ7").Value = 200350000
